$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 250
$ws.Range("I6").Value = 250
$ws.Range("K6").Value = 750
$ws.Range("M6").Value = -638

$ws.Range("H28").Value = 968.1429000000001
$ws.Range("J28").Value = 1233
$ws.Range("L28").Value = 1233
$ws.Range("N28").Value = -2203

$ws.Range("H33").Value = 22790.223
$ws.Range("I33").Value = 716
$ws.Range("K33").Value = 716
$ws.Range("M33").Value = -487

$ws.Range("H53").Value = 81
$ws.Range("I53").Value = 71.5
$ws.Range("K53").Value = 71.5
$ws.Range("M53").Value = 565.5

$ws.Range("H98").Value = 1133.5
$ws.Range("I98").Value = 1152.7778
$ws.Range("K98").Value = 1152.7778
$ws.Range("M98").Value = 345.2221999999999

$ws.Range("H108").Value = 150000
$ws.Range("J108").Value = 150000
$ws.Range("L108").Value = 150000
$ws.Range("N108").Value = -157680

$ws.Range("H111").Value = 950
$ws.Range("I111").Value = 950
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2850
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 217
$ws.Range("N111").ClearContents()

$ws.Range("H122").Value = 1133.5
$ws.Range("I122").Value = 1152.7778
$ws.Range("K122").Value = 3458.3334
$ws.Range("M122").Value = -1008.3334

$ws.Range("H137").Value = 600
$ws.Range("I137").Value = 600
$ws.Range("K137").Value = 1800
$ws.Range("M137").Value = 750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1514.4572
$ws.Range("I32").Value = 1559.2188
$ws.Range("J32").Value = 1037
$ws.Range("K32").Value = 1559.2188
$ws.Range("L32").Value = 1037
$ws.Range("M32").Value = -1272.2188
$ws.Range("N32").Value = -1611

$ws.Range("H61").Value = 3250
$ws.Range("I61").Value = 1400
$ws.Range("J61").Value = 4571.4287
$ws.Range("K61").Value = 1400
$ws.Range("L61").Value = 4571.4287
$ws.Range("M61").Value = -1188
$ws.Range("N61").Value = -4995.4287

$ws.Range("H122").Value = 2454.3635
$ws.Range("I122").Value = 1499.75
$ws.Range("K122").Value = 4499.25
$ws.Range("M122").Value = -2049.25

$ws.Range("H132").Value = 2190.0688
$ws.Range("I132").Value = 1119.8096
$ws.Range("K132").Value = 3359.4288
$ws.Range("M132").Value = -829.4288000000001

$ws.Range("H136").Value = 3250
$ws.Range("I136").Value = 1400
$ws.Range("J136").Value = 4571.4287
$ws.Range("K136").Value = 4200
$ws.Range("L136").Value = 13714.2861
$ws.Range("M136").Value = -1650
$ws.Range("N136").Value = -18814.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2304.5715
$ws.Range("I31").Value = 1225.1765
$ws.Range("K31").Value = 1225.1765
$ws.Range("M31").Value = -930.1765

$ws.Range("H34").Value = 2304.5715
$ws.Range("I34").Value = 1225.1765
$ws.Range("K34").Value = 1225.1765
$ws.Range("M34").Value = -1023.1765

$ws.Range("H58").Value = 4567.3335
$ws.Range("I58").Value = 3401
$ws.Range("K58").Value = 3401
$ws.Range("M58").Value = -3198

$ws.Range("H132").Value = 3749.5
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws.Range("H136").Value = 4567.3335
$ws.Range("I136").Value = 3401
$ws.Range("K136").Value = 10203
$ws.Range("M136").Value = -7653

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 937.5
$ws.Range("J98").Value = 916.6667
$ws.Range("L98").Value = 2750.0001
$ws.Range("N98").Value = -5746.0001

$ws.Range("H105").Value = 9499.5
$ws.Range("I105").Value = 9000
$ws.Range("J105").Value = 9999
$ws.Range("K105").Value = 27000
$ws.Range("L105").Value = 29997
$ws.Range("M105").Value = -24379
$ws.Range("N105").Value = -35239

$ws.Range("H109").Value = 3716.9092
$ws.Range("I109").Value = 1814.3334
$ws.Range("K109").Value = 5443.0002
$ws.Range("M109").Value = -4403.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9999.5
$ws.Range("I70").Value = 9999.5
$ws.Range("K70").Value = 9999.5
$ws.Range("M70").Value = -9729.5

$ws.Range("H73").Value = 9999.5
$ws.Range("I73").Value = 9999.5
$ws.Range("K73").Value = 9999.5
$ws.Range("M73").Value = -9063.5

$ws.Range("H103").Value = 47500
$ws.Range("J103").Value = 47500
$ws.Range("L103").Value = 47500
$ws.Range("N103").Value = -49844

$ws.Range("H113").Value = 2957.8572
$ws.Range("I113").Value = 2927
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 2927
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = -757
$ws.Range("N113").Value = -7339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1104.4615
$ws.Range("I22").Value = 775
$ws.Range("J22").Value = 1250.8889
$ws.Range("K22").Value = 775
$ws.Range("L22").Value = 1250.8889
$ws.Range("M22").Value = -480
$ws.Range("N22").Value = -1840.8889

$ws.Range("H27").Value = 1104.4615
$ws.Range("I27").Value = 775
$ws.Range("J27").Value = 1250.8889
$ws.Range("K27").Value = 775
$ws.Range("L27").Value = 1250.8889
$ws.Range("M27").Value = -668
$ws.Range("N27").Value = -1464.8889

$ws.Range("H61").Value = 998.6667
$ws.Range("I61").Value = 998
$ws.Range("K61").Value = 998
$ws.Range("M61").Value = -796

$ws.Range("H113").Value = 998.6667
$ws.Range("I113").Value = 998
$ws.Range("K113").Value = 998
$ws.Range("M113").Value = 1172

$ws.Range("H115").Value = 29800
$ws.Range("J115").Value = 29800
$ws.Range("L115").Value = 29800
$ws.Range("N115").Value = -32150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2578.2856
$ws.Range("I122").Value = 2409.6
$ws.Range("K122").Value = 7228.799999999999
$ws.Range("M122").Value = -4778.799999999999

$ws.Range("H132").Value = 3003.7083
$ws.Range("I132").Value = 1818.1875
$ws.Range("J132").Value = 5374.75
$ws.Range("K132").Value = 5454.5625
$ws.Range("L132").Value = 15998.7495
$ws.Range("M132").Value = -2924.5625
$ws.Range("N132").Value = -21184.25
